$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("E2").Value = "Handed back: in sync with en-US"
$ovw.Range("F2").Value = "Handed back: in sync with en-US"
$ovw.Columns.Item(5).ColumnWidth = 29.16666666666667
$ovw.Columns.Item(6).ColumnWidth = 29.16666666666667

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("K2").Value = "2016-08-28 06:54:29"
$zhcn.Range("P2").Value = ""
$zhcn.Columns.Item(3).ColumnWidth = 29.16666666666667
$zhcn.Columns.Item(16).ColumnWidth = 12.833333333333334

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("K2").Value = "2016-08-28 06:54:36"
$dede.Range("P2").Value = ""
$dede.Columns.Item(3).ColumnWidth = 29.16666666666667
$dede.Columns.Item(16).ColumnWidth = 12.833333333333334
